$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '76.468.64'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.049.73'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +4.39%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '202.24'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '625.46'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +4.52%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.209'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +5.97%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '3.048.83'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +4.38%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.437'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.14'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +5.10%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.613.36'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +4.44%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '29.45'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +4.99%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '76.433.52'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000194'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.055.83'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +4.48%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.52'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +4.68%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '9.06'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +3.53%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '374.90'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.32'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.35'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.24%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '73.58'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +2.89%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.205.83'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +4.43%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.43'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +4.37%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +2.35%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.48%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.30'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +7.15%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '509.00'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +6.64%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '20.87'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +3.22%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.27%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.386'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +6.47%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '20.03'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '190.74'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +4.33%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.113'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.806'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +22.80%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +4.07%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.27'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +6.72%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '42.26'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +5.59%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.68'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.47'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +3.80%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +6.71%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.90'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +4.72%  '

Write-Host "Applied 92 cell updates"
